$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "39.484.16"
Set-TextValue "E2" "  +2.01%  "

Set-TextValue "D3" "2.163.60"
Set-TextValue "E3" "  +3.07%  "

Set-TextValue "E4" "  -0.01%  "

Set-TextValue "D5" "227.85"
Set-TextValue "E5" "  -0.28%  "

Set-TextValue "E6" "  +1.13%  "

Set-TextValue "D7" "64.23"
Set-TextValue "E7" "  +3.70%  "

Set-TextValue "E8" "  +0.01%  "

Set-TextValue "E9" "  +2.48%  "

Set-TextValue "E10" "  +2.04%  "

Set-TextValue "E11" "  +0.10%  "

Set-TextValue "D12" "16.27"
Set-TextValue "E12" "  +2.65%  "

Set-TextValue "D13" "2.484.77"
Set-TextValue "E13" "  +3.04%  "

Set-TextValue "D14" "22.16"
Set-TextValue "E14" "  +0.70%  "

Set-TextValue "D15" "0.816"
Set-TextValue "E15" "  +1.41%  "

Set-TextValue "E16" "  +0.46%  "

Set-TextValue "D17" "2.164.74"
Set-TextValue "E17" "  +3.71%  "

Set-TextValue "D18" "39.448.26"
Set-TextValue "E18" "  +1.85%  "

Set-TextValue "E19" "  +0.20%  "

Set-TextValue "D20" "6.14"
Set-TextValue "E20" "  +1.28%  "

Set-TextValue "D21" "0.0₃0853"
Set-TextValue "E21" "  +1.71%  "

Set-TextValue "D22" "229.73"
Set-TextValue "E22" "  +0.92%  "

Set-TextValue "E23" "  +0.05%  "

Set-TextValue "B24" "Toncoin"
Set-TextValue "C24" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D24" "2.35"
Set-TextValue "E24" "  +0.07%  "

Set-TextValue "B25" "PancakeSwap"
Set-TextValue "C25" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D25" "2.36"
Set-TextValue "E25" "  +1.40%  "

Set-TextValue "D26" "9.61"
Set-TextValue "E26" "  +0.55%  "

Set-TextValue "D27" "172.17"
Set-TextValue "E27" "  -0.03%  "

Set-TextValue "D28" "0.141"
Set-TextValue "E28" "  +1.74%  "

Set-TextValue "D29" "1.46"
Set-TextValue "E29" "  +2.36%  "

Set-TextValue "D30" "19.93"
Set-TextValue "E30" "  +3.20%  "

Set-TextValue "E31" "  +2.98%  "

Set-TextValue "D33" "4.62"
Set-TextValue "E33" "  +1.72%  "

Set-TextValue "D34" "7.12"
Set-TextValue "E34" "  +5.13%  "

Set-TextValue "D35" "4.72"
Set-TextValue "E35" "  -0.66%  "

Set-TextValue "D36" "0.0618"
Set-TextValue "E36" "  -0.19%  "

Set-TextValue "D37" "2.45"
Set-TextValue "E37" "  +1.45%  "

Set-TextValue "E38" "  +0.13%  "

Set-TextValue "E39" "  +0.12%  "

Set-TextValue "B40" "Aave"
Set-TextValue "C40" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D40" "103.44"
Set-TextValue "E40" "  +0.95%  "

Set-TextValue "B41" "VeChain"
Set-TextValue "C41" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D41" "0.0230"
Set-TextValue "E41" "  +0.87%  "

Set-TextValue "D42" "17.86"
Set-TextValue "E42" "  -1.60%  "

Set-TextValue "D43" "1.526.72"
Set-TextValue "E43" "  -0.46%  "

Set-TextValue "E44" "  +3.88%  "

Set-TextValue "D45" "0.0933"
Set-TextValue "E45" "  +2.61%  "

Set-TextValue "D46" "2.83"
Set-TextValue "E46" "  +0.77%  "

Set-TextValue "D47" "4.29"
Set-TextValue "E47" "  +4.31%  "

Set-TextValue "E48" "  +5.78%  "

Set-TextValue "D49" "7.77"
Set-TextValue "E49" "  -0.69%  "

Set-TextValue "D50" "2.367.93"
Set-TextValue "E50" "  +3.26%  "

Set-TextValue "E51" "  -0.35%  "
